# Append run: 2025-12-11 12:40 JST
# Re-writes the "ランサーズ" listing sheet with the latest scrape:
#  - all existing rows get the new scrape timestamp
#  - 3 brand-new listings are inserted, the whole list re-sorted by
#    priority score (column G) descending, which shuffles several rows
#  - column B/D get a little wider to fit the new (longer) text

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ts = "2025-12-11 12:40:37"

# --- row data, already sorted by priority score (G) descending ------------
# Each entry: Title, Category, Price, Deadline, Url, Score, Skills(optional)
$rows = @(
    @{B="【急募】AIチャットボット開発のプロフェッショナルを探しています!"; C="システム開発"; D="100,000 円 ~ 200,000 円 / 固定"; E="期限情報なし"; F="https://www.lancers.jp/work/detail/5451734"; G=368; H="🔥AI,Ai ◆開発"},
    @{B="【募集】RPAツール「RoboTANGO」設定代行の専門家を探しています"; C="システム開発"; D="50,000 円 ~ 100,000 円 / 固定"; E="期限情報なし"; F="https://www.lancers.jp/work/detail/5405023"; G=178; H="★bot ◆ツール"},
    @{B="Webシステム開発"; C="システム開発"; D="200,000 円 ~ 300,000 円 / 固定"; E="期限情報なし"; F="https://www.lancers.jp/work/detail/5451859"; G=118; H="◆開発,システム開発"},
    @{B="自動出品システムの開発"; C="システム開発"; D="100,000 円 ~ 200,000 円 / 固定"; E="期限情報なし"; F="https://www.lancers.jp/work/detail/5451514"; G=83; H="◆開発"},
    @{B="【愛知県近辺 在住の方希望 / リモート相談可能】経験豊富なWebフロントエンド開発エンジニア募集!"; C="システム開発"; D="500,000 円 ~ 1,000,000 円 / 固定"; E="期限情報なし"; F="https://www.lancers.jp/work/detail/5451972"; G=75; H="◆開発"},
    @{B="Access DB家賃管理SYSを最新Access で稼働できるように"; C="システム開発"; D="100,000 円 ~ 200,000 円 / 固定"; E="期限情報なし"; F="https://www.lancers.jp/work/detail/5451626"; G=38; H="◇管理"},
    @{B="進行管理およびチームディレクションを担当"; C="システム開発"; D="~ 5,000 円 / 固定"; E="期限情報なし"; F="https://www.lancers.jp/work/detail/5418064"; G=30; H="◇管理"},
    @{B="Rubyの暗号化機能のPHP化"; C="システム開発"; D="20,000 円 ~ 50,000 円 / 固定"; E="期限情報なし"; F="https://www.lancers.jp/work/detail/5451714"; G=28; H="○PHP"},
    @{B="【オンライン講師募集】バックエンドの基礎を教えていただける方"; C="システム開発"; D="100,000 円 ~ 200,000 円 / 固定"; E="期限情報なし"; F="https://www.lancers.jp/work/detail/5451420"; G=18; H=$null},
    @{B="注目 限定公開 PR 限定公開の仕事"; C="システム開発"; D="20,000 円 ~ 50,000 円 / 固定"; E="期限情報なし"; F="https://www.lancers.jp/work/detail/5450323"; G=13; H=$null},
    @{B="【急募】desknetスタンダード版からNeo版への移行サポート"; C="システム開発"; D="20,000 円 ~ 50,000 円 / 固定"; E="期限情報なし"; F="https://www.lancers.jp/work/detail/5451838"; G=13; H=$null}
)

# --- write the cell values -------------------------------------------------
$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $ts
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
    $ws.Cells.Item($r, 7).Value = $row.G
    if ($row.H) {
        $ws.Cells.Item($r, 8).Value = $row.H
    }
    $r = $r + 1
}

# --- hyperlinks: drop the old set, re-add one per row, in row order -------
$ws.Hyperlinks.Delete()
$r = 2
foreach ($row in $rows) {
    $ws.Hyperlinks.Add($ws.Cells.Item($r, 6), $row.F)
    $r = $r + 1
}

# --- column widths (B: 38 -> 52, D: 28 -> 30) ------------------------------
$ws.Columns.Item(2).ColumnWidth = 52 - 11/12 + 0.03
$ws.Columns.Item(4).ColumnWidth = 30 - 11/12 + 0.03
